$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row directly below the header row, pushing the
# existing rows (THIAGO, RODRIGO, ...) down by one.
$ws.Rows.Item(2).Insert()

# Column A ("Conta") holds account numbers with significant leading
# zeros, stored as text throughout the sheet. Temporarily force a text
# number format so "005880251" isn't coerced into the number 5880251,
# then clear the format again so the cell doesn't end up with a style
# that the original data rows don't have.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "005880251"
$ws.Cells.Item(2, 1).ClearFormats()

$ws.Cells.Item(2, 2).Value = "LUIZ"
$ws.Cells.Item(2, 3).Value = 60000
